$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Red font color used to mark rows/cells whose columns were deleted in a
# later data-collection wave ("day 8"). This reuses the same red font
# already used elsewhere in the sheet for "deleted day 6" markers.
$red = 255

# --- Row 3 (gender): mark deleted, add new "deleted day 8" note in E3 ---
$ws.Range("E3").Value = "deleted day 8"
$ws.Range("A3:D3").Font.Color = $red
$ws.Range("E3").Font.Color = $red
$ws.Range("G3").Font.Color = $red

# --- Row 12 (base_temp_txt): update E12 note, mark deleted ---
$ws.Range("E12").Value = "new, deleted day 8"
$ws.Range("A12:B12").Font.Color = $red
$ws.Range("D12").Font.Color = $red
$ws.Range("E12").Font.Color = $red
$ws.Range("G12").Font.Color = $red

# --- Row 15 (base_esr_txt): update E15 note, mark deleted ---
$ws.Range("E15").Value = "new, deleted day 8"
$ws.Range("A15:B15").Font.Color = $red
$ws.Range("D15").Font.Color = $red
$ws.Range("E15").Font.Color = $red
$ws.Range("G15").Font.Color = $red

# --- Row 17 (base_cavitation_txt): update E17 note, mark deleted ---
$ws.Range("E17").Value = "new, deleted day 8"
$ws.Range("A17:B17").Font.Color = $red
$ws.Range("D17").Font.Color = $red
$ws.Range("E17").Font.Color = $red
$ws.Range("G17").Font.Color = $red

# --- Row 20 (strep_resistance_txt): update E20 note, mark deleted ---
$ws.Range("E20").Value = "new, deleted day 8"
$ws.Range("A20").Font.Color = $red
$ws.Range("D20").Font.Color = $red
$ws.Range("E20").Font.Color = $red
$ws.Range("G20").Font.Color = $red

# --- Row 21 (strep_resistance_range): update E21 note, mark deleted ---
$ws.Range("E21").Value = "new, deleted day 8"
$ws.Range("A21").Font.Color = $red
$ws.Range("D21").Font.Color = $red
$ws.Range("E21").Font.Color = $red
$ws.Range("G21").Font.Color = $red

# --- Row 24 (radiologic_6mon_txt): already "new, deleted day 6"; mark deleted ---
$ws.Range("A24:B24").Font.Color = $red
$ws.Range("D24").Font.Color = $red
$ws.Range("E24").Font.Color = $red
$ws.Range("G24").Font.Color = $red

# --- Row 25 (rad_num): add "deleted day 8" note in E25, mark deleted ---
$ws.Range("E25").Value = "deleted day 8"
$ws.Range("A25:C25").Font.Color = $red
$ws.Range("D25").Font.Color = $red
$ws.Range("E25").Font.Color = $red
$ws.Range("G25").Font.Color = $red

# --- Row 26 (improved): add "deleted day 8" note in E26, mark deleted ---
$ws.Range("E26").Value = "deleted day 8"
$ws.Range("A26:C26").Font.Color = $red
$ws.Range("D26").Font.Color = $red
$ws.Range("E26").Font.Color = $red
$ws.Range("G26").Font.Color = $red

# Match the author's final cursor position in the saved file.
[void]$ws.Range("E27").Select()
